$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row appended at the bottom of the data range (row 55), mirroring the
# formatting of the previous last row (row 54: bold/bordered index cell in
# column A, date-time numeric format in column E).
$ws.Range("A54:V54").Copy()
$ws.Range("A55:V55").PasteSpecial(-4122)

$row = 55
$ws.Cells.Item($row, 1).Value = 54
$ws.Cells.Item($row, 2).Value = "france"
$ws.Cells.Item($row, 3).Value = "ligue-1"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45195.875
$ws.Cells.Item($row, 6).Value = "Lille"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = "Reims"
$ws.Cells.Item($row, 9).Value = 2
$ws.Cells.Item($row, 10).Value = 1.76
$ws.Cells.Item($row, 11).Value = "05/09/2023 12:01"
$ws.Cells.Item($row, 12).Value = 1.9
$ws.Cells.Item($row, 13).Value = "26/09/2023 20:56"
$ws.Cells.Item($row, 14).Value = 4.09
$ws.Cells.Item($row, 15).Value = "05/09/2023 12:01"
$ws.Cells.Item($row, 16).Value = 3.93
$ws.Cells.Item($row, 17).Value = "26/09/2023 20:58"
$ws.Cells.Item($row, 18).Value = 4.47
$ws.Cells.Item($row, 19).Value = "05/09/2023 12:01"
$ws.Cells.Item($row, 20).Value = 4.13
$ws.Cells.Item($row, 21).Value = "26/09/2023 20:57"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/france/ligue-1/lille-reims/nNmvumia/"
